$d = $word.ActiveDocument

$d.Content.Find.Execute("Logging", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Login", 2)
